$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.201.60"
$ws.Range("E2").Value = "  -3.37%  "
$ws.Range("D3").Value = "'1.609.42"
$ws.Range("E3").Value = "  -2.71%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'0.9996"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "'302.85"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("D7").Value = "'0.3780"
$ws.Range("E7").Value = "  -3.26%  "
$ws.Range("D8").Value = "'0.3674"
$ws.Range("E8").Value = "  -4.50%  "
$ws.Range("E9").Value = "  -4.16%  "
$ws.Range("D10").Value = "'0.9994"
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").Value = "'1.278"
$ws.Range("E11").Value = "  -5.94%  "
$ws.Range("D12").Value = "'0.08105"
$ws.Range("E12").Value = "  -4.17%  "
$ws.Range("D13").Value = "'23.15"
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("D14").Value = "'6.628"
$ws.Range("E14").Value = "  -7.05%  "
$ws.Range("D15").Value = "'7.631"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").Value = "'0.00001272"
$ws.Range("E16").Value = "  -3.49%  "
$ws.Range("D17").Value = "'1.609.15"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "'91.57"
$ws.Range("E18").Value = "  -3.32%  "
$ws.Range("D19").Value = "'0.06795"
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("D20").Value = "'18.43"
$ws.Range("E20").Value = "  -7.07%  "
$ws.Range("D21").Value = "'6.602"
$ws.Range("E21").Value = "  -4.77%  "
$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  -4.29%  "
$ws.Range("D24").Value = "'23.210.98"
$ws.Range("E24").Value = "  -3.32%  "
$ws.Range("D25").Value = "'2.362"
$ws.Range("E25").Value = "  -5.00%  "
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("D27").Value = "'21.12"
$ws.Range("E27").Value = "  -4.62%  "
$ws.Range("D28").Value = "'150.50"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").Value = "'5.256"
$ws.Range("E29").Value = "  -3.60%  "
$ws.Range("D30").Value = "'132.60"
$ws.Range("E30").Value = "  -4.94%  "
$ws.Range("D31").Value = "'2.417"
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("D32").Value = "'7.018"
$ws.Range("E32").Value = "  -11.15%  "
$ws.Range("D33").Value = "'1.787.29"
$ws.Range("E33").Value = "  -2.66%  "
$ws.Range("D34").Value = "'0.9895"
$ws.Range("E34").Value = "  -5.35%  "
$ws.Range("D35").Value = "'0.07754"
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("D36").Value = "'0.02791"
$ws.Range("E36").Value = "  -6.27%  "
$ws.Range("D37").Value = "'6.324"
$ws.Range("E37").Value = "  -6.66%  "
$ws.Range("D38").Value = "'0.2555"
$ws.Range("E38").Value = "  -4.96%  "
$ws.Range("D39").Value = "'10.14"
$ws.Range("E39").Value = "  -7.03%  "
$ws.Range("D40").Value = "'0.08878"
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("D41").Value = "'1.397"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D42").Value = "'0.7190"
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("D43").Value = "'12.82"
$ws.Range("E43").Value = "  -4.86%  "
$ws.Range("D44").Value = "'15.92"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("D45").Value = "'0.6637"
$ws.Range("E45").Value = "  -4.70%  "
$ws.Range("D46").Value = "'2.309"
$ws.Range("E46").Value = "  -6.15%  "
$ws.Range("D47").Value = "'0.9986"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "'3.979"
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("D49").Value = "'0.08020"
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("D50").Value = "'132.02"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("E51").Value = "  -3.40%  "
